$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Range("E26").Value = 392.5
$ws.Range("E27").Value = 552.617
$ws.Range("E28").Value = 385.424
$ws.Range("E29").Value = 390.18
